# Kennedy Gichuhi Ngaruiya Q0510 - "adding averages and more checks"
#
# 1) Training Dashboard: PERIOD TO EXPIRE (col H) drops by 8 days and
#    LAST UPDATE (col I) moves from 08-Sep-2025 -> 16-Sep-2025 for rows 3-12.
# 2) Exam Dashboard: widen the COMMENTS-adjacent column (col E, "MARKS ATTAINED"
#    width) from 10 -> 15 chars, and reword the pass remark from "OK" to
#    "date is valid" for rows 3-5.
# 3) Header / title styling: bold white text on the title row and the
#    column-header row (both sheets), title font size reset to the default 11.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Training Dashboard
$ws2 = $wb.Worksheets.Item(2)   # Exam Dashboard

# ---------------------------------------------------------------------------
# 1) Training Dashboard data updates (H3:H12, I3:I12)
# ---------------------------------------------------------------------------
$periodToExpire = @{
    3  = 334
    4  = 638
    5  = 364
    6  = 448
    7  = 352
    8  = 423
    9  = 588
    10 = 677
    11 = 86
    12 = 181
}

foreach ($row in 3..12) {
    $ws1.Range("H$row").Value = $periodToExpire[$row]

    $dateCell = $ws1.Range("I$row")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 2) Exam Dashboard: column width + remark text
# ---------------------------------------------------------------------------
$ws2.Range("E1").EntireColumn.ColumnWidth = 14.17

foreach ($row in 3..5) {
    $ws2.Range("E$row").Value = "date is valid"
}

# ---------------------------------------------------------------------------
# 3) Header / title font styling (both sheets) -> bold white
# ---------------------------------------------------------------------------
$white = 16777215

# Column-header rows first (already bold, size 11) - only color changes.
$ws1.Range("A2:K2").Font.Color = $white
$ws2.Range("A2:G2").Font.Color = $white

# Title cells: reset size 14 -> 11, then apply the same white color.
$ws1.Range("A1").Font.Size  = 11
$ws1.Range("A1").Font.Color = $white

$ws2.Range("A1").Font.Size  = 11
$ws2.Range("A1").Font.Color = $white
